$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# B4 holds the Cypher query used to populate the "FilesTab" export.
# Drop the `File Type` and `Breed` columns from the RETURN clause.
$ws.Range("B4").Value = "MATCH (f:file)-->(parent)`nWITH DISTINCT f, parent`nMATCH (f)-[*]->(c:case)<--(demo:demographic)`n MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`n MATCH (samp:sample)-->(c) `nWHERE samp.summarized_sample_type IN [`"Primary Malignant Tumor Tissue`"] `nWITH DISTINCT f, parent, c, demo, diag, s`nRETURN coalesce(f.file_name, '') AS ``File Name``, `n        coalesce(labels(parent)[0], '') AS ``Association``,`n        coalesce(f.file_description, '') AS ``Description``,`n        coalesce(f.file_format, '') AS ``Format``,`n        coalesce(f.file_size, '') AS ``Size``,`n        coalesce(c.case_id, '') AS ``Case ID``, `n        coalesce(diag.disease_term,'') AS Diagnosis , `n        coalesce(s.clinical_study_designation,'') AS ``Study Code``"

# Active selection in the sheet moves from D4 to B4.
$ws.Range("B4").Select()
